# Features commit for booking flight
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the booking reference values in column A (rows 2-5)
$ws.Range("A2").Value = 9363339066
$ws.Range("A3").Value = 9363339066
$ws.Range("A4").Value = 9363339066
$ws.Range("A5").Value = 9363339066

# Move the active selection to A6 (also resets the scrolled view back to column A)
$ws.Range("A6").Select()
